# Updated cryptos list on Sun Dec  3 05:08:15 UTC 2023 with GitHub Actions
# Refresh the coin price / 1h-volume-change table with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "39.382.62"
$ws.Cells.Item(2, 5).Value = "  +1.51%  "

$ws.Cells.Item(3, 4).Value = "2.156.54"
$ws.Cells.Item(3, 5).Value = "  +3.02%  "

$ws.Cells.Item(4, 5).Value = "  +0.06%  "

$ws.Cells.Item(5, 4).Value = "'227.51"
$ws.Cells.Item(5, 5).Value = "  -0.61%  "

$ws.Cells.Item(6, 5).Value = "  +0.88%  "

$ws.Cells.Item(7, 4).Value = "'64.00"
$ws.Cells.Item(7, 5).Value = "  +3.98%  "

$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(9, 5).Value = "  +2.77%  "

$ws.Cells.Item(10, 4).Value = "'0.0858"
$ws.Cells.Item(10, 5).Value = "  +1.47%  "

$ws.Cells.Item(11, 5).Value = "  -0.08%  "

$ws.Cells.Item(12, 4).Value = "'15.96"
$ws.Cells.Item(12, 5).Value = "  +4.24%  "

$ws.Cells.Item(13, 4).Value = "2.476.74"
$ws.Cells.Item(13, 5).Value = "  +2.97%  "

$ws.Cells.Item(14, 4).Value = "'22.23"
$ws.Cells.Item(14, 5).Value = "  +0.66%  "

$ws.Cells.Item(15, 4).Value = "'0.809"
$ws.Cells.Item(15, 5).Value = "  +0.45%  "

$ws.Cells.Item(16, 5).Value = "  +1.09%  "

$ws.Cells.Item(17, 4).Value = "2.160.98"
$ws.Cells.Item(17, 5).Value = "  +3.37%  "

$ws.Cells.Item(18, 4).Value = "39.334.49"
$ws.Cells.Item(18, 5).Value = "  +1.56%  "

$ws.Cells.Item(19, 4).Value = "'71.76"
$ws.Cells.Item(19, 5).Value = "  -0.23%  "

$ws.Cells.Item(20, 4).Value = "'6.11"
$ws.Cells.Item(20, 5).Value = "  +0.60%  "

$ws.Cells.Item(21, 5).Value = "  +1.54%  "

$ws.Cells.Item(22, 4).Value = "'231.30"
$ws.Cells.Item(22, 5).Value = "  +1.40%  "

$ws.Cells.Item(23, 5).Value = "  +0.00%  "

$ws.Cells.Item(24, 4).Value = "'2.44"
$ws.Cells.Item(24, 5).Value = "  +2.95%  "

$ws.Cells.Item(25, 4).Value = "'2.35"
$ws.Cells.Item(25, 5).Value = "  +0.54%  "

$ws.Cells.Item(26, 2).Value = "Cosmos"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(26, 4).Value = "'9.66"
$ws.Cells.Item(26, 5).Value = "  +1.28%  "

$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(27, 4).Value = "'172.06"
$ws.Cells.Item(27, 5).Value = "  +0.32%  "

$ws.Cells.Item(28, 5).Value = "  +1.16%  "

$ws.Cells.Item(29, 4).Value = "'19.88"
$ws.Cells.Item(29, 5).Value = "  +2.98%  "

$ws.Cells.Item(30, 4).Value = "'1.40"
$ws.Cells.Item(30, 5).Value = "  -1.14%  "

$ws.Cells.Item(31, 5).Value = "  +9.41%  "

$ws.Cells.Item(32, 5).Value = "  +0.75%  "

$ws.Cells.Item(33, 5).Value = "  +2.16%  "

$ws.Cells.Item(34, 4).Value = "'4.74"
$ws.Cells.Item(34, 5).Value = "  -0.04%  "

$ws.Cells.Item(35, 5).Value = "  +8.72%  "

$ws.Cells.Item(36, 4).Value = "'0.0617"
$ws.Cells.Item(36, 5).Value = "  -0.34%  "

$ws.Cells.Item(37, 5).Value = "  +0.60%  "

$ws.Cells.Item(38, 4).Value = "'3.59"
$ws.Cells.Item(38, 5).Value = "  -0.34%  "

$ws.Cells.Item(39, 5).Value = "  +0.05%  "

$ws.Cells.Item(40, 4).Value = "'103.83"
$ws.Cells.Item(40, 5).Value = "  +2.63%  "

$ws.Cells.Item(42, 4).Value = "'17.79"
$ws.Cells.Item(42, 5).Value = "  -2.58%  "

$ws.Cells.Item(43, 4).Value = "1.538.92"
$ws.Cells.Item(43, 5).Value = "  +0.31%  "

$ws.Cells.Item(44, 5).Value = "  +3.59%  "

$ws.Cells.Item(45, 4).Value = "'7.92"
$ws.Cells.Item(45, 5).Value = "  +3.59%  "

$ws.Cells.Item(46, 5).Value = "  +0.47%  "

$ws.Cells.Item(47, 4).Value = "'0.0924"
$ws.Cells.Item(47, 5).Value = "  +1.45%  "

$ws.Cells.Item(48, 5).Value = "  +5.62%  "

$ws.Cells.Item(49, 4).Value = "'4.22"
$ws.Cells.Item(49, 5).Value = "  +2.92%  "

$ws.Cells.Item(50, 4).Value = "2.360.39"
$ws.Cells.Item(50, 5).Value = "  +3.02%  "

$ws.Cells.Item(51, 5).Value = "  +0.03%  "
